$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.111.08'
$ws.Range("E2").Value = '  -1.20%  '

$ws.Range("D3").Value = '1.790.44'
$ws.Range("E3").Value = '  -1.61%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.54'
$ws.Range("E5").Value = '  -0.44%  '

$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5189'
$ws.Range("E7").Value = '  +1.91%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3805'
$ws.Range("E8").Value = '  -3.68%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07969'
$ws.Range("E9").Value = '  -4.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.44'
$ws.Range("E10").Value = '  -0.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.091'
$ws.Range("E11").Value = '  -1.75%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.262'
$ws.Range("E12").Value = '  -0.95%  '

$ws.Range("B13").Value = 'BinanceUSD'
$ws.Range("C13").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.003'
$ws.Range("E13").Value = '  +0.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.48'
$ws.Range("E14").Value = '  -2.73%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.271'
$ws.Range("E15").Value = '  -3.38%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.787.35'
$ws.Range("E16").Value = '  -1.60%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.12'
$ws.Range("E17").Value = '  -1.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001090'
$ws.Range("E18").Value = '  -5.00%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06545'
$ws.Range("E19").Value = '  -1.62%  '

$ws.Range("E20").Value = '  +0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.27'
$ws.Range("E21").Value = '  -2.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.941'
$ws.Range("E22").Value = '  -2.98%  '

$ws.Range("D23").Value = '28.138.21'
$ws.Range("E23").Value = '  -1.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.11'
$ws.Range("E24").Value = '  -3.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.259'
$ws.Range("E25").Value = '  -0.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.01'
$ws.Range("E26").Value = '  +2.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.37'
$ws.Range("E27").Value = '  -4.31%  '

$ws.Range("D28").Value = '1.993.00'
$ws.Range("E28").Value = '  -1.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.331'
$ws.Range("E29").Value = '  -3.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.66'
$ws.Range("E30").Value = '  -2.39%  '

$ws.Range("E31").Value = '  -0.93%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.053'
$ws.Range("E32").Value = '  -5.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.675'
$ws.Range("E33").Value = '  +0.70%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.526'
$ws.Range("E34").Value = '  -4.59%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07193'
$ws.Range("E35").Value = '  +1.79%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.17'
$ws.Range("E36").Value = '  +7.69%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02308'
$ws.Range("E37").Value = '  -1.27%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2138'
$ws.Range("E38").Value = '  -3.91%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.074'
$ws.Range("E39").Value = '  -2.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.563'
$ws.Range("E40").Value = '  -3.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6152'
$ws.Range("E41").Value = '  -2.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.158'
$ws.Range("E42").Value = '  -1.60%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.366'
$ws.Range("E43").Value = '  -2.48%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.15'
$ws.Range("E44").Value = '  -2.82%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.759'
$ws.Range("E45").Value = '  +0.74%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5919'
$ws.Range("E46").Value = '  -0.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '126.85'
$ws.Range("E47").Value = '  +1.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.213'
$ws.Range("E48").Value = '  +2.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.912'
$ws.Range("E49").Value = '  -3.61%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06756'
$ws.Range("E50").Value = '  -1.96%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.40'
$ws.Range("E51").Value = '  -2.62%  '
